$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: bump the date serial by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# D32:D35 - updated unit prices
$ws.Range("D32").Value = 5775
$ws.Range("D33").Value = 4587
$ws.Range("D34").Value = 4576
$ws.Range("D35").Value = 6583.5
